$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (date, volume, prices) per diff ---
$ws.Range("D2").Value = 44235
$ws.Range("M2").Value = 15
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 25000
$ws.Range("S2").Value = 1250
$ws.Range("D3").Value = 44432
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 24000
$ws.Range("S3").Value = 1200
$ws.Range("D4").Value = 44428
$ws.Range("D5").Value = 44454
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 25000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 25000
$ws.Range("S5").Value = 1250
$ws.Range("D6").Value = 44424
$ws.Range("M6").Value = 25
$ws.Range("D7").Value = 44398
$ws.Range("M7").Value = 15
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 25000
$ws.Range("S7").Value = 1250
$ws.Range("D8").Value = 44421
$ws.Range("M8").Value = 20
$ws.Range("D9").Value = 44222
$ws.Range("D10").Value = 44349
$ws.Range("D11").Value = 44396
$ws.Range("M11").Value = 12
$ws.Range("D12").Value = 44412
$ws.Range("M12").Value = 20
$ws.Range("D13").Value = 44431
$ws.Range("M13").Value = 40
$ws.Range("D14").Value = 44232
$ws.Range("M14").Value = 15
$ws.Range("N14").Value = 25000
$ws.Range("O14").Value = 25000
$ws.Range("P14").Value = 25000
$ws.Range("S14").Value = 1250
$ws.Range("D15").Value = 44231
$ws.Range("M15").Value = 15
$ws.Range("D16").Value = 44435
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 24000
$ws.Range("O16").Value = 24000
$ws.Range("P16").Value = 24000
$ws.Range("S16").Value = 1200
$ws.Range("D17").Value = 44238
$ws.Range("M17").Value = 30
$ws.Range("D18").Value = 44391
$ws.Range("D19").Value = 44334
$ws.Range("M19").Value = 20
$ws.Range("D20").Value = 44419
$ws.Range("M20").Value = 40
$ws.Range("D21").Value = 44434
$ws.Range("M21").Value = 20
$ws.Range("N21").Value = 24000
$ws.Range("O21").Value = 24000
$ws.Range("P21").Value = 24000
$ws.Range("S21").Value = 1200
$ws.Range("D22").Value = 44452
$ws.Range("M22").Value = 25
$ws.Range("N22").Value = 25000
$ws.Range("O22").Value = 25000
$ws.Range("P22").Value = 25000
$ws.Range("S22").Value = 1250
$ws.Range("D24").Value = 44400
$ws.Range("M24").Value = 5
$ws.Range("D25").Value = 44392
$ws.Range("M25").Value = 10
$ws.Range("D26").Value = 44442
$ws.Range("M26").Value = 25
$ws.Range("N26").Value = 23000
$ws.Range("O26").Value = 23000
$ws.Range("P26").Value = 23000
$ws.Range("S26").Value = 1150
$ws.Range("D27").Value = 44214
$ws.Range("D29").Value = 44418
$ws.Range("M29").Value = 20
$ws.Range("N29").Value = 24000
$ws.Range("O29").Value = 24000
$ws.Range("P29").Value = 24000
$ws.Range("S29").Value = 1200
$ws.Range("D30").Value = 44414
$ws.Range("M30").Value = 15
$ws.Range("N30").Value = 25000
$ws.Range("O30").Value = 25000
$ws.Range("P30").Value = 25000
$ws.Range("S30").Value = 1250
$ws.Range("D31").Value = 44389
$ws.Range("M31").Value = 20
$ws.Range("D32").Value = 44249
$ws.Range("D33").Value = 44390
$ws.Range("M33").Value = 10
$ws.Range("D34").Value = 44251
$ws.Range("M34").Value = 15
$ws.Range("N34").Value = 25000
$ws.Range("O34").Value = 25000
$ws.Range("P34").Value = 25000
$ws.Range("S34").Value = 1250
$ws.Range("D35").Value = 44433
$ws.Range("M35").Value = 10
$ws.Range("D36").Value = 44221
$ws.Range("M36").Value = 30
$ws.Range("N36").Value = 25000
$ws.Range("O36").Value = 25000
$ws.Range("P36").Value = 25000
$ws.Range("S36").Value = 1250
$ws.Range("D37").Value = 44363
$ws.Range("M37").Value = 30
$ws.Range("N37").Value = 24000
$ws.Range("O37").Value = 24000
$ws.Range("P37").Value = 24000
$ws.Range("S37").Value = 1200
$ws.Range("D38").Value = 44356
$ws.Range("M38").Value = 15
$ws.Range("N38").Value = 24000
$ws.Range("O38").Value = 24000
$ws.Range("P38").Value = 24000
$ws.Range("S38").Value = 1200
$ws.Range("D39").Value = 44175
$ws.Range("M39").Value = 25
$ws.Range("N39").Value = 23000
$ws.Range("O39").Value = 23000
$ws.Range("P39").Value = 23000
$ws.Range("S39").Value = 1150

# --- Append new row 40 ---
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "Vega Modelo de Temuco"
$ws.Range("C40").Value = "La Araucanía"
$ws.Range("D40").Value = 44425
$ws.Range("D40").NumberFormat = $ws.Range("D39").NumberFormat
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100108
$ws.Range("H40").Value = "Tropicales y subtropicales"
$ws.Range("I40").Value = 100108007
$ws.Range("J40").Value = "Coco"
$ws.Range("K40").Value = "Sin especificar"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 15
$ws.Range("N40").Value = 24000
$ws.Range("O40").Value = 24000
$ws.Range("P40").Value = 24000
$ws.Range("Q40").Value = "$/malla 20 unidades"
$ws.Range("R40").Value = "Perú"
$ws.Range("S40").Value = 1200
$ws.Range("T40").Value = 20
